# Change the table style applied to every table in the deck from the
# default "Table_0" style ({0F348C7D-CC2B-498C-98C9-12C344C2484E}) to the
# built-in PowerPoint table style {5A98252D-CE5C-461E-A920-4279631CC404}
# (this mirrors picking a different swatch in the Table Styles gallery
# while a table is selected).

$p = $ppt.ActivePresentation
$newStyleId = "{5A98252D-CE5C-461E-A920-4279631CC404}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
